$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("ED2B")
$ws3 = $wb.Worksheets.Item("ED2C")

# ED2B (sheet2) updates to column A
$ws2.Cells.Item(2, 1).Value = 49.33110367892976
$ws2.Cells.Item(3, 1).Value = 10.45751633986928
$ws2.Cells.Item(4, 1).Value = 11.11111111111111
$ws2.Cells.Item(5, 1).Value = 11.43790849673203
$ws2.Cells.Item(6, 1).Value = 10.7843137254902
$ws2.Cells.Item(7, 1).Value = 10.99656357388316
$ws2.Cells.Item(8, 1).Value = 9.965635738831615
$ws2.Cells.Item(9, 1).Value = 14.26116838487972
$ws2.Cells.Item(10, 1).Value = 11.85567010309278
$ws2.Cells.Item(12, 1).Value = 8.524590163934427
$ws2.Cells.Item(13, 1).Value = 10.81967213114754
$ws2.Cells.Item(14, 1).Value = 11.47540983606557
$ws2.Cells.Item(15, 1).Value = 12.45901639344262
$ws2.Cells.Item(16, 1).Value = 29.07563025210084
$ws2.Cells.Item(17, 1).Value = 47.1571906354515
$ws2.Cells.Item(18, 1).Value = 48.82943143812709
$ws2.Cells.Item(19, 1).Value = 14.70588235294118
$ws2.Cells.Item(20, 1).Value = 11.76470588235294
$ws2.Cells.Item(21, 1).Value = 12.02749140893471
$ws2.Cells.Item(22, 1).Value = 11.51202749140894
$ws2.Cells.Item(23, 1).Value = 45.48494983277592
$ws2.Cells.Item(24, 1).Value = 50.33444816053512
$ws2.Cells.Item(25, 1).Value = 12.41830065359477
$ws2.Cells.Item(26, 1).Value = 12.09150326797386
$ws2.Cells.Item(27, 1).Value = 12.37113402061856
$ws2.Cells.Item(28, 1).Value = 9.965635738831615
$ws2.Cells.Item(29, 1).Value = 49.49832775919732

# ED2C (sheet3) updates to column A
$ws3.Cells.Item(2, 1).Value = 49.33110367892976
$ws3.Cells.Item(3, 1).Value = 10.45751633986928
$ws3.Cells.Item(4, 1).Value = 11.43790849673203
$ws3.Cells.Item(5, 1).Value = 11.43790849673203
$ws3.Cells.Item(6, 1).Value = 10.7843137254902
$ws3.Cells.Item(8, 1).Value = 44.31438127090301
$ws3.Cells.Item(9, 1).Value = 46.32107023411371
$ws3.Cells.Item(10, 1).Value = 10.99656357388316
$ws3.Cells.Item(11, 1).Value = 9.965635738831615
$ws3.Cells.Item(12, 1).Value = 14.26116838487972
$ws3.Cells.Item(13, 1).Value = 11.51202749140894
$ws3.Cells.Item(14, 1).Value = 13.39869281045752
$ws3.Cells.Item(15, 1).Value = 13.39869281045752
$ws3.Cells.Item(16, 1).Value = 8.196721311475411
$ws3.Cells.Item(17, 1).Value = 10.16393442622951
$ws3.Cells.Item(18, 1).Value = 11.80327868852459
$ws3.Cells.Item(19, 1).Value = 10.81967213114754
$ws3.Cells.Item(20, 1).Value = 14.43298969072165
$ws3.Cells.Item(21, 1).Value = 12.88659793814433
$ws3.Cells.Item(22, 1).Value = 29.41176470588236
